$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    3  = -4
    4  = 1
    5  = 0
    6  = 4
    7  = 2
    8  = -2
    9  = 1
    10 = 5
    11 = 2
    12 = -3
    13 = 5
    14 = 3
    17 = 1
    18 = 5
    19 = -5
    20 = -5
    21 = -3
    23 = -4
    24 = 7
    25 = -10
    26 = -5
    27 = -9
    29 = -2
    30 = 0
    31 = -1
    32 = -6
    33 = 0
    34 = 2
    35 = 1
    36 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
